$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 780, pushing existing rows 780-830 down to 782-832
$rng = $ws.Range("A780:T781")
$rng.Insert()

# Fill new row 780 with fresh data
$ws.Cells.Item(780, 1).Value = 5
$ws.Cells.Item(780, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(780, 3).Value = "Maule"
$ws.Cells.Item(780, 4).Value = 45021
$ws.Cells.Item(780, 5).Value = 7
$ws.Cells.Item(780, 6).Value = "Fruta"
$ws.Cells.Item(780, 7).Value = 100102
$ws.Cells.Item(780, 8).Value = "Cítricos"
$ws.Cells.Item(780, 9).Value = 100102005
$ws.Cells.Item(780, 10).Value = "Naranja"
$ws.Cells.Item(780, 11).Value = "Valencia"
$ws.Cells.Item(780, 12).Value = "Primera"
$ws.Cells.Item(780, 13).Value = 200
$ws.Cells.Item(780, 14).Value = 14000
$ws.Cells.Item(780, 15).Value = 14000
$ws.Cells.Item(780, 16).Value = 14000
$ws.Cells.Item(780, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(780, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(780, 19).Value = 933
$ws.Cells.Item(780, 20).Value = 15

# Fill new row 781 with fresh data
$ws.Cells.Item(781, 1).Value = 5
$ws.Cells.Item(781, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(781, 3).Value = "Maule"
$ws.Cells.Item(781, 4).Value = 45021
$ws.Cells.Item(781, 5).Value = 7
$ws.Cells.Item(781, 6).Value = "Fruta"
$ws.Cells.Item(781, 7).Value = 100102
$ws.Cells.Item(781, 8).Value = "Cítricos"
$ws.Cells.Item(781, 9).Value = 100102005
$ws.Cells.Item(781, 10).Value = "Naranja"
$ws.Cells.Item(781, 11).Value = "Valencia"
$ws.Cells.Item(781, 12).Value = "Primera"
$ws.Cells.Item(781, 13).Value = 360
$ws.Cells.Item(781, 14).Value = 13000
$ws.Cells.Item(781, 15).Value = 13000
$ws.Cells.Item(781, 16).Value = 13000
$ws.Cells.Item(781, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(781, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(781, 19).Value = 867
$ws.Cells.Item(781, 20).Value = 15
